# Apply updated loading_percent values for Case_5_15 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new value
$newValues = @{
    2 = @{ "B" = 16.97360316259524; "C" = 8.927199038821646; "D" = 8.179484323685447; "E" = 12.71991215803186; "F" = 34.81600840799211; "H" = 7.344005520526261; "I" = 25.90020916364016; "J" = 9.934216296803267; "L" = 11.31073906751209; "M" = 16.43699828793618; "O" = 26.75469176677645 }
    3 = @{ "B" = 16.47134532421446; "C" = 8.610807771909339; "D" = 8.181699412178087; "E" = 12.75325384499563; "F" = 34.94969786146713; "H" = 7.344005520526261; "I" = 26.04598948743862; "J" = 9.952701352754108; "L" = 11.30793128581069; "M" = 16.31609476631338; "O" = 26.87018873949428 }
    4 = @{ "B" = 16.15609732352274; "C" = 8.409388972004313; "D" = 8.183768361102507; "E" = 12.77484751703678; "F" = 35.04026403197653; "H" = 7.344005520526261; "I" = 26.14085487560537; "J" = 9.964656550594725; "L" = 11.30733878506446; "M" = 16.24291824750621; "O" = 26.94731461301878 }
    5 = @{ "B" = 16.02608497996058; "C" = 8.325597113346845; "D" = 8.184790320582319; "E" = 12.7839299496335; "F" = 35.07929759271922; "H" = 7.344005520526261; "I" = 26.18086078127814; "J" = 9.969681047319192; "L" = 11.30738307239281; "M" = 16.21338629803003; "O" = 26.98030246486561 }
    6 = @{ "B" = 16.00440862042997; "C" = 8.311582857480191; "D" = 8.18497083492357; "E" = 12.78545518576075; "F" = 35.08590739916269; "H" = 7.344005520526261; "I" = 26.18758514009088; "J" = 9.970524594479176; "L" = 11.30740771814485; "M" = 16.20850056857323; "O" = 26.98587410431334 }
    7 = @{ "B" = 16.15434994644205; "C" = 8.408265736152352; "D" = 8.183781418689765; "E" = 12.77496885968676; "F" = 35.04078184681683; "H" = 7.344005520526261; "I" = 26.14138895226683; "J" = 9.964723694018025; "L" = 11.30733822385021; "M" = 16.24251877439453; "O" = 26.94775319294547 }
    8 = @{ "B" = 16.80195149638629; "C" = 8.819641057639853; "D" = 8.180101320800279; "E" = 12.73117601561069; "F" = 34.86034128753486; "H" = 7.344005520526261; "I" = 25.94936309808275; "J" = 9.940464585190183; "L" = 11.30953685048111; "M" = 16.39510425348552; "O" = 26.79322427583579 }
    9 = @{ "B" = 18.0102943951129; "C" = 9.566254641836032; "D" = 8.178482941893904; "E" = 12.65416392544987; "F" = 34.574017428469; "H" = 7.344005520526261; "I" = 25.61526115507183; "J" = 9.897674795464003; "L" = 11.32277126833792; "M" = 16.7017474894972; "O" = 26.53961049476535 }
    10 = @{ "B" = 18.85185638328612; "C" = 10.0743084607989; "D" = 8.18067102501773; "E" = 12.60293877286647; "F" = 34.40512485986784; "H" = 7.344005520526261; "I" = 25.39563289959223; "J" = 9.869123757873906; "L" = 11.33785814413247; "M" = 16.93023627193091; "O" = 26.38359764035931 }
    11 = @{ "B" = 19.22310808565301; "C" = 10.29599013201311; "D" = 8.182391654653731; "E" = 12.58078763064008; "F" = 34.33735817567422; "H" = 7.344005520526261; "I" = 25.30131964981396; "J" = 9.85675584473349; "L" = 11.34586855992479; "M" = 17.03458952136333; "O" = 26.31924689785462 }
    12 = @{ "B" = 19.36191317263303; "C" = 10.37853314289658; "D" = 8.183146759054209; "E" = 12.57256436227206; "F" = 34.31300482435236; "H" = 7.344005520526261; "I" = 25.26641014116439; "J" = 9.852161147404974; "L" = 11.349065228991; "M" = 17.07414040991737; "O" = 26.29583428928692 }
    13 = @{ "B" = 19.3320998748416; "C" = 10.36081910070722; "D" = 8.182979541155037; "E" = 12.57432806765886; "F" = 34.31819148303418; "H" = 7.344005520526261; "I" = 25.27389272952698; "J" = 9.85314675635545; "L" = 11.34836953401606; "M" = 17.06562127386704; "O" = 26.30083405467169 }
    14 = @{ "B" = 19.23456395004575; "C" = 10.30280935732248; "D" = 8.182451706874676; "E" = 12.58010779670312; "F" = 34.33532835965313; "H" = 7.344005520526261; "I" = 25.29843149198691; "J" = 9.856376059590247; "L" = 11.34612828899235; "M" = 17.03784284529539; "O" = 26.31730155577261 }
    15 = @{ "B" = 19.17458538522665; "C" = 10.26709272997464; "D" = 8.182141854202431; "E" = 12.5836695007009; "F" = 34.34599573417494; "H" = 7.344005520526261; "I" = 25.31356701417621; "J" = 9.858365648632441; "L" = 11.34477667666296; "M" = 17.02083152533735; "O" = 26.32751293265038 }
    16 = @{ "B" = 18.82735111915035; "C" = 10.05962708846336; "D" = 8.180573108257818; "E" = 12.60440951310393; "F" = 34.40973639354107; "H" = 7.344005520526261; "I" = 25.40190911677256; "J" = 9.869944473206923; "L" = 11.33735758938091; "M" = 16.92342272677951; "O" = 26.38793664144499 }
    17 = @{ "B" = 18.61128290216378; "C" = 9.929902215621924; "D" = 8.179795963706656; "E" = 12.61742725881988; "F" = 34.45116429197963; "H" = 7.344005520526261; "I" = 25.45753761538365; "J" = 9.877206238496395; "L" = 11.33309890169629; "M" = 16.86375362990942; "O" = 26.42670284410622 }
    18 = @{ "B" = 18.48592237212869; "C" = 9.85440222063443; "D" = 8.179417282404563; "E" = 12.62502314116985; "F" = 34.47584532152922; "H" = 7.344005520526261; "I" = 25.49006043894939; "J" = 9.881441404328873; "L" = 11.33075751410527; "M" = 16.82947385590327; "O" = 26.4496230571139 }
    19 = @{ "B" = 18.44329500360845; "C" = 9.8286886196592; "D" = 8.17930082246429; "E" = 12.62761361857056; "F" = 34.48434820316517; "H" = 7.344005520526261; "I" = 25.50116258304361; "J" = 9.882885399779985; "L" = 11.32998337631866; "M" = 16.81787499185853; "O" = 26.45749032282473 }
    20 = @{ "B" = 18.63439677267414; "C" = 9.943803669946623; "D" = 8.179871627432551; "E" = 12.6160302815692; "F" = 34.44666592221579; "H" = 7.344005520526261; "I" = 25.45156134650795; "J" = 9.876427170939632; "L" = 11.33354107060516; "M" = 16.87010151438101; "O" = 26.42251161942849 }
    21 = @{ "B" = 19.26326177629173; "C" = 10.31988664142158; "D" = 8.182603940706548; "E" = 12.57840568059617; "F" = 34.33025929383927; "H" = 7.344005520526261; "I" = 25.29120201764239; "J" = 9.855425129167958; "L" = 11.34678217891999; "M" = 17.04600130159352; "O" = 26.31243868755565 }
    22 = @{ "B" = 19.66384199829103; "C" = 10.55748345235344; "D" = 8.18499280260799; "E" = 12.55477661080576; "F" = 34.26180951221469; "H" = 7.344005520526261; "I" = 25.19108943265715; "J" = 9.84221627347342; "L" = 11.35638700127515; "M" = 17.16115211202992; "O" = 26.24607140104292 }
    23 = @{ "B" = 19.45103303406447; "C" = 10.43143715227801; "D" = 8.183662890795775; "E" = 12.56730020367159; "F" = 34.29764279744926; "H" = 7.344005520526261; "I" = 25.24409208464454; "J" = 9.849218898403267; "L" = 11.35117428738685; "M" = 17.09968471494021; "O" = 26.28098187044399 }
    24 = @{ "B" = 18.62395053392397; "C" = 9.937521682752342; "D" = 8.179837207664814; "E" = 12.61666150666631; "F" = 34.44869694556186; "H" = 7.344005520526261; "I" = 25.45426153165138; "J" = 9.87677919957971; "L" = 11.33334083265834; "M" = 16.86723155671622; "O" = 26.42440450057791 }
    25 = @{ "B" = 17.69094625406605; "C" = 9.371138443177427; "D" = 8.178325005762662; "E" = 12.67405371194164; "F" = 34.64421464938716; "H" = 7.344005520526261; "I" = 25.70110424617196; "J" = 9.908741582628435; "L" = 11.31824330138010; "M" = 16.61813367952541; "O" = 26.60290943722659 }
}

foreach ($rowNum in $newValues.Keys) {
    $rowMap = $newValues[$rowNum]
    foreach ($colLetter in $rowMap.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $rowMap[$colLetter]
    }
}
